$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: Grand Manan - corrected figures
$ws.Range("B38").Value = 103165300
$ws.Range("C38").Value = 552000
$ws.Range("D38").Value = 5635600
$ws.Range("E38").Value = 109352900
$ws.Range("G38").Value = 866300
$ws.Range("I38").Value = 17579300
$ws.Range("J38").Value = 126932200
$ws.Range("K38").Value = 135721850
$ws.Range("L38").Value = 135682764

# Row 39: previously-duplicated "Grand Manan" row is actually Grande-Anse
$ws.Range("A39").Value = "Grande-Anse"
$ws.Range("B39").Value = 20258600
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 3606700
$ws.Range("E39").Value = 23865300
$ws.Range("G39").Value = 82300
$ws.Range("I39").Value = 4817300
$ws.Range("J39").Value = 28682600
$ws.Range("K39").Value = 31091250
$ws.Range("L39").Value = 31070375

# Rows 40-104: municipality names shift down by one position because of the inserted Grande-Anse row
$ws.Range("A40").Value = "Hampton"
$ws.Range("A41").Value = "Hartland"
$ws.Range("A42").Value = "Hillsborough"
$ws.Range("A43").Value = "Kedgwick"
$ws.Range("A44").Value = "Lamèque"
$ws.Range("A45").Value = "McAdam"
$ws.Range("A46").Value = "Memramcook"
$ws.Range("A47").Value = "Minto"
$ws.Range("A48").Value = "Nackawic"
$ws.Range("A49").Value = "Néguac"
$ws.Range("A50").Value = "Perth-Andover"
$ws.Range("A51").Value = "Petit-Rocher"
$ws.Range("A52").Value = "Petitcodiac"
$ws.Range("A53").Value = "Plaster Rock"
$ws.Range("A54").Value = "Rexton"
$ws.Range("A55").Value = "Richibucto"
$ws.Range("A56").Value = "Rogersville"
$ws.Range("A57").Value = "Saint Andrews"
$ws.Range("A58").Value = "Saint George"
$ws.Range("A59").Value = "Saint-Antoine"
$ws.Range("A60").Value = "Saint-Louis-de-Kent"
$ws.Range("A61").Value = "Saint-Quentin"
$ws.Range("A62").Value = "Salisbury"
$ws.Range("A63").Value = "Saint-Léonard"
$ws.Range("A64").Value = "Alma"
$ws.Range("A65").Value = "Aroostook"
$ws.Range("A66").Value = "Baker-Brook"
$ws.Range("A67").Value = "Balmoral"
$ws.Range("A68").Value = "Bas-Caraquet"
$ws.Range("A69").Value = "Bath"
$ws.Range("A70").Value = "Bertrand"
$ws.Range("A71").Value = "Blackville"
$ws.Range("A72").Value = "Bristol"
$ws.Range("A73").Value = "Cambridge-Narrows"
$ws.Range("A74").Value = "Canterbury"
$ws.Range("A75").Value = "Centreville"
$ws.Range("A76").Value = "Dorchester"
$ws.Range("A77").Value = "Drummond"
$ws.Range("A78").Value = "Fredericton Junction"
$ws.Range("A79").Value = "Gagetown"
$ws.Range("A80").Value = "Harvey"
$ws.Range("A81").Value = "Le Goulet"
$ws.Range("A82").Value = "Maisonnette"
$ws.Range("A83").Value = "Meductic"
$ws.Range("A84").Value = "Millville"
$ws.Range("A85").Value = "Nigadoo"
$ws.Range("A86").Value = "Norton"
$ws.Range("A87").Value = "Paquetville"
$ws.Range("A88").Value = "Pointe-Verte"
$ws.Range("A89").Value = "Port Elgin"
$ws.Range("A90").Value = "Riverside-Albert"
$ws.Range("A91").Value = "Rivière-Verte"
$ws.Range("A92").Value = "Saint-François-de-Madawaska"
$ws.Range("A93").Value = "Saint-Isidore"
$ws.Range("A94").Value = "Saint-Léolin"
$ws.Range("A95").Value = "Sainte-Anne-de-Madawaska"
$ws.Range("A96").Value = "Sainte-Marie-Saint-Raphaël"
$ws.Range("A97").Value = "Saint-André"
$ws.Range("A98").Value = "Saint-Hilaire"
$ws.Range("A99").Value = "St. Martins"
$ws.Range("A100").Value = "Stanley"
$ws.Range("A101").Value = "Sussex Corner"
$ws.Range("A102").Value = "Tide Head"
$ws.Range("A103").Value = "Tracy"
$ws.Range("A104").Value = "Lac Baker"
